$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''278.19'
$ws.Range('E2').Value = '''6.58%'
$ws.Range('D3').Value = '''27.40'
$ws.Range('E3').Value = '''1.99%'
$ws.Range('D4').Value = '''4.841'
$ws.Range('E4').Value = '''3.10%'
$ws.Range('D5').Value = '''0.06250'
$ws.Range('E5').Value = '''0.56%'
$ws.Range('D6').Value = '''6.900'
$ws.Range('E6').Value = '''2.13%'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').Value = '''3.269'
$ws.Range('E7').Value = '''2.80%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = '''0.8786'
$ws.Range('E8').Value = '''3.20%'
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D9').Value = '''0.9458'
$ws.Range('E9').Value = '''3.65%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '''0.1449'
$ws.Range('E10').Value = '''3.24%'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').Value = '''0.05222'
$ws.Range('E11').Value = '''6.24%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '''0.07293'
$ws.Range('E12').Value = '''2.83%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '''0.03160'
$ws.Range('E13').Value = '''1.94%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '''0.09056'
$ws.Range('E14').Value = '''0.11%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '''0.001552'
$ws.Range('E15').Value = '''1.27%'
$ws.Range('B16').Value = 'One'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D16').Value = '''0.0006267'
$ws.Range('E16').Value = '''1.78%'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').Value = '''0.006099'
$ws.Range('E17').Value = '''1.04%'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').Value = '''3.460'
$ws.Range('E18').Value = '''0.45%'
$ws.Range('D19').Value = '''2.285'
$ws.Range('E19').Value = '''6.53%'
$ws.Range('E20').Value = '''-0.62%'
$ws.Range('D21').Value = '''0.1310'
$ws.Range('E21').Value = '''-0.02%'
$ws.Range('D22').Value = '''3.850'
$ws.Range('E22').Value = '''-5.81%'
$ws.Range('D23').Value = '''0.04326'
$ws.Range('E23').Value = '''2.23%'
$ws.Range('D24').Value = '''0.001175'
$ws.Range('E24').Value = '''-2.43%'
$ws.Range('D25').Value = '''0.004275'
$ws.Range('E25').Value = '''4.85%'
$ws.Range('D26').Value = '''0.0001200'
$ws.Range('E26').Value = '''-0.04%'
$ws.Range('E27').Value = '''3.08%'
$ws.Range('D40').Value = '''0.04027'
$ws.Range('E40').Value = '''1.82%'
$ws.Range('D41').Value = '''0.006411'
$ws.Range('E41').Value = '''55.06%'
$ws.Range('D42').Value = '''0.1154'
$ws.Range('E42').Value = '''3.77%'
$ws.Range('D43').Value = '''0.002109'
$ws.Range('E43').Value = '''-4.61%'
$ws.Range('D44').Value = '''0.01200'
$ws.Range('E44').Value = '''-13.60%'
$ws.Range('D45').Value = '''0.00005087'
$ws.Range('E45').Value = '''-1.45%'
$ws.Range('E46').Value = '''-0.03%'
$ws.Range('D47').Value = '''2.375'
$ws.Range('E47').Value = '''813.73%'
$ws.Range('D49').Value = '''0.00002100'
$ws.Range('E49').Value = '''-0.03%'
$ws.Range('D50').Value = '''0.0002000'
$ws.Range('E50').Value = '''-0.03%'
